# chore: adapt column header formatting to respective input file names (#7)
#
# The sheet holds a side-by-side AHB diff: the left half of the header row
# used a generic "_old" suffix and the right half a generic "_new" suffix.
# Rename them to the concrete format-version names ("_FV2310" / "_FV2404"),
# turn the used range into a proper Excel Table (so the header row gets an
# AutoFilter), and freeze the header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header cells: "_old" -> "_FV2310" (columns A..J) ------------
$ws.Cells.Item(1, 1).Value = "Segmentname_FV2310"
$ws.Cells.Item(1, 2).Value = "Segmentgruppe_FV2310"
$ws.Cells.Item(1, 3).Value = "Segment_FV2310"
$ws.Cells.Item(1, 4).Value = "Datenelement_FV2310"
$ws.Cells.Item(1, 5).Value = "Segment ID_FV2310"
$ws.Cells.Item(1, 6).Value = "Code_FV2310"
$ws.Cells.Item(1, 7).Value = "Qualifier_FV2310"
$ws.Cells.Item(1, 8).Value = "Beschreibung_FV2310"
$ws.Cells.Item(1, 9).Value = "Bedingungsausdruck_FV2310"
$ws.Cells.Item(1, 10).Value = "Bedingung_FV2310"

# Column K ("diff") stays as-is.

# --- "_new" -> "_FV2404" (columns L..U) -------------------------------------
$ws.Cells.Item(1, 12).Value = "Segmentname_FV2404"
$ws.Cells.Item(1, 13).Value = "Segmentgruppe_FV2404"
$ws.Cells.Item(1, 14).Value = "Segment_FV2404"
$ws.Cells.Item(1, 15).Value = "Datenelement_FV2404"
$ws.Cells.Item(1, 16).Value = "Segment ID_FV2404"
$ws.Cells.Item(1, 17).Value = "Code_FV2404"
$ws.Cells.Item(1, 18).Value = "Qualifier_FV2404"
$ws.Cells.Item(1, 19).Value = "Beschreibung_FV2404"
$ws.Cells.Item(1, 20).Value = "Bedingungsausdruck_FV2404"
$ws.Cells.Item(1, 21).Value = "Bedingung_FV2404"

# --- 2. Turn the used range into an Excel Table with an AutoFilter ---------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U72"), $null, 1)
$tbl.Name = "Table1"

# --- 3. Freeze the header row ------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
